# Generate Report for Handback
# Swap the "36038c19-..." and "fbebf677-..." rows (row 2 <-> row 3) across the
# Overview / zh-cn / de-de sheets, and update their status/dates to reflect a
# completed handback.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws.Range("B2").Value = "e2e\36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws.Range("G2").Value = "2016-08-23 00:45:57"

$ws.Range("A3").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws.Range("B3").Value = "e2e\fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-08-23 00:44:50"

$ws.Hyperlinks.Item(1).TextToDisplay = "e2e\36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws.Hyperlinks.Item(2).TextToDisplay = "e2e\fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws.Range("G2").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.4ad38bf24ece5dede3f5c797292c36779823b837.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-23 00:45:52"
$ws.Range("I2").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws.Range("J2").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.4ad38bf24ece5dede3f5c797292c36779823b837.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-23 00:46:14"

$ws.Range("A3").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.3180788ccb1d5c858ef0ef8e59d53f0fa210ab48.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-23 00:44:45"
$ws.Range("I3").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws.Range("J3").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.3180788ccb1d5c858ef0ef8e59d53f0fa210ab48.zh-cn.xlf"
$ws.Range("P3").Value = ""

$ws.Hyperlinks.Item(1).TextToDisplay = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws.Hyperlinks.Item(2).TextToDisplay = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws.Hyperlinks.Item(3).TextToDisplay = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws.Hyperlinks.Item(4).TextToDisplay = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"

$ws.Columns.Item(16).ColumnWidth = 13.7470528738839

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws.Range("G2").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.4ad38bf24ece5dede3f5c797292c36779823b837.de-de.xlf"
$ws.Range("H2").Value = "2016-08-23 00:45:57"
$ws.Range("I2").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws.Range("J2").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.4ad38bf24ece5dede3f5c797292c36779823b837.de-de.xlf"
$ws.Range("K2").Value = "2016-08-23 00:46:21"

$ws.Range("A3").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.3180788ccb1d5c858ef0ef8e59d53f0fa210ab48.de-de.xlf"
$ws.Range("H3").Value = "2016-08-23 00:44:50"
$ws.Range("I3").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws.Range("J3").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.3180788ccb1d5c858ef0ef8e59d53f0fa210ab48.de-de.xlf"
$ws.Range("P3").Value = ""

$ws.Hyperlinks.Item(1).TextToDisplay = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws.Hyperlinks.Item(2).TextToDisplay = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws.Hyperlinks.Item(3).TextToDisplay = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws.Hyperlinks.Item(4).TextToDisplay = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"

$ws.Columns.Item(16).ColumnWidth = 13.7470528738839
